$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recover dropped data: the "ID Competição" column (B) was mistakenly
# truncated to 67 for every entrant row; restore the correct value 267
# for each data row (B2:B73).
$ws.Range("B2:B73").Value = 267
